$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.65", "0.584") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.537.91'
$ws.Range('E2').Value = '  -1.05%  '

$ws.Range('D3').Value = '2.604.94'
$ws.Range('E3').Value = '  -0.84%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '589.70'
$ws.Range('E5').Value = '  -2.13%  '

$ws.Range('D6').Value = '149.43'
$ws.Range('E6').Value = '  -1.01%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '0.584'
$ws.Range('E8').Value = '  -1.01%  '

$ws.Range('E9').Value = '  -0.49%  '

$ws.Range('D10').Value = '5.75'
$ws.Range('E10').Value = '  +0.28%  '

$ws.Range('D11').Value = '0.387'
$ws.Range('E11').Value = '  +0.29%  '

$ws.Range('E12').Value = '  +0.30%  '

$ws.Range('D13').Value = '27.63'
$ws.Range('E13').Value = '  -0.15%  '

$ws.Range('D14').Value = '3.070.65'
$ws.Range('E14').Value = '  -0.94%  '

$ws.Range('D15').Value = '63.354.26'
$ws.Range('E15').Value = '  -1.10%  '

$ws.Range('D16').Value = '0.0000156'
$ws.Range('E16').Value = '  +4.54%  '

$ws.Range('D17').Value = '2.583.49'
$ws.Range('E17').Value = '  -1.58%  '

$ws.Range('D18').Value = '12.08'
$ws.Range('E18').Value = '  -0.95%  '

$ws.Range('D19').Value = '4.73'
$ws.Range('E19').Value = '  +1.57%  '

$ws.Range('D20').Value = '345.23'
$ws.Range('E20').Value = '  -1.96%  '

$ws.Range('D21').Value = '6.86'
$ws.Range('E21').Value = '  -1.75%  '

$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').Value = '66.78'
$ws.Range('E23').Value = '  +0.16%  '

$ws.Range('D24').Value = '1.69'
$ws.Range('E24').Value = '  -3.86%  '

$ws.Range('D25').Value = '9.23'
$ws.Range('E25').Value = '  -0.15%  '

$ws.Range('D26').Value = '1.65'
$ws.Range('E26').Value = '  -2.68%  '

$ws.Range('D27').Value = '8.30'
$ws.Range('E27').Value = '  +2.56%  '

$ws.Range('D28').Value = '549.83'
$ws.Range('E28').Value = '  +2.28%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '0.161'
$ws.Range('E29').Value = '  -2.63%  '

$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.07%  '

$ws.Range('D31').Value = '2.03'
$ws.Range('E31').Value = '  -1.77%  '

$ws.Range('D32').Value = '0.0₃0860'
$ws.Range('E32').Value = '  +0.52%  '

$ws.Range('D33').Value = '1.77'
$ws.Range('E33').Value = '  +1.32%  '

$ws.Range('D34').Value = '5.31'
$ws.Range('E34').Value = '  +0.47%  '

$ws.Range('D35').Value = '6.08'
$ws.Range('E35').Value = '  -1.06%  '

$ws.Range('D36').Value = '166.04'
$ws.Range('E36').Value = '  -1.06%  '

$ws.Range('D37').Value = '0.412'
$ws.Range('E37').Value = '  +0.31%  '

$ws.Range('E38').Value = '  -0.07%  '

$ws.Range('D39').Value = '19.43'
$ws.Range('E39').Value = '  -0.57%  '

$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  -3.48%  '

$ws.Range('E41').Value = '  -0.05%  '

$ws.Range('D42').Value = '164.98'
$ws.Range('E42').Value = '  -3.03%  '

$ws.Range('D43').Value = '4.03'
$ws.Range('E43').Value = '  +2.40%  '

$ws.Range('D44').Value = '22.93'
$ws.Range('E44').Value = '  +6.35%  '

$ws.Range('D45').Value = '0.0580'
$ws.Range('E45').Value = '  -1.73%  '

$ws.Range('D46').Value = '2.10'
$ws.Range('E46').Value = '  +5.04%  '

$ws.Range('D47').Value = '0.632'
$ws.Range('E47').Value = '  +0.46%  '

$ws.Range('D48').Value = '0.0250'
$ws.Range('E48').Value = '  +1.27%  '

$ws.Range('D49').Value = '0.0959'
$ws.Range('E49').Value = '  -0.91%  '

$ws.Range('D50').Value = '19.09'
$ws.Range('E50').Value = '  -1.15%  '

$ws.Range('D51').Value = '0.0₆0228'
$ws.Range('E51').Value = '  +14.90%  '

# Restore default cell style (clears the temporary text-number-format)
# while keeping the values stored as text.
$ws.Range("D2:E51").Style = "Normal"
